$wb = $excel.ActiveWorkbook

# "Repayment Schedule" is the 4th sheet (sheet4.xml / rId4).
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14) on the Repayment Schedule
# sheet. This shifts the old N -> O, O -> P, P -> Q, matching the diff
# (dimension A1:P12 -> A1:Q12, header "Late" N1 -> O1, "Outstanding" P1 -> Q1).
$wsSchedule.Columns.Item(14).Insert()

# Update the selection on the Repayment Schedule sheet to R4, and make it
# the active/selected tab.
$wsSchedule.Range("R4").Select() | Out-Null
$wsSchedule.Activate()
